$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pretest change: SNAP2 chassis # for the first antenna (row 4) 22 -> 21 ---
$ws.Range("K4").Value = 21

# --- Column width adjustments (every column narrowed by roughly one pixel) ---
# Columns A, J and O keep their original/default width; every other column
# gets a slightly smaller width. The ColumnWidth values below are chosen so
# that, after the host engine's internal rounding, the serialized widths in
# the saved workbook land as close as possible to the target widths.
$ws.Range("B1").EntireColumn.ColumnWidth = 5.5
$ws.Range("C1").EntireColumn.ColumnWidth = 8.333333333333334
$ws.Range("D1").EntireColumn.ColumnWidth = 9.333333333333334
$ws.Range("E1").EntireColumn.ColumnWidth = 8.166666666666666
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 5.333333333333333
$ws.Range("H1").EntireColumn.ColumnWidth = 6.833333333333333
$ws.Range("I1").EntireColumn.ColumnWidth = 14.666666666666666
$ws.Range("K1").EntireColumn.ColumnWidth = 5.833333333333333
$ws.Range("L1:M1").EntireColumn.ColumnWidth = 6.333333333333333
$ws.Range("N1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("P1").EntireColumn.ColumnWidth = 11.666666666666666
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.333333333333333
$ws.Range("R1:S1").EntireColumn.ColumnWidth = 6.666666666666667
$ws.Range("T1").EntireColumn.ColumnWidth = 5.166666666666667
$ws.Range("U1:V1").EntireColumn.ColumnWidth = 5.666666666666667
$ws.Range("W1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("X1").EntireColumn.ColumnWidth = 6.833333333333333
$ws.Range("Y1:Z1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AA1").EntireColumn.ColumnWidth = 24.833333333333332
